$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.582.49"
$ws.Range("E2").Value = "'  -1.39%  "
$ws.Range("D3").Value = "'3.089.67"
$ws.Range("E3").Value = "'  -2.00%  "
$ws.Range("E4").Value = "'  -0.17%  "
$ws.Range("D5").Value = "'588.85"
$ws.Range("E5").Value = "'  -0.93%  "
$ws.Range("D6").Value = "'153.25"
$ws.Range("E6").Value = "'  +4.28%  "
$ws.Range("E7").Value = "'  -0.06%  "
$ws.Range("D8").Value = "'0.554"
$ws.Range("E8").Value = "'  +4.21%  "
$ws.Range("D9").Value = "'3.087.13"
$ws.Range("E9").Value = "'  -2.02%  "
$ws.Range("E10").Value = "'  -2.12%  "
$ws.Range("E11").Value = "'  -0.29%  "
$ws.Range("D12").Value = "'0.462"
$ws.Range("E12").Value = "'  -0.69%  "
$ws.Range("D13").Value = "'37.70"
$ws.Range("E13").Value = "'  +0.68%  "
$ws.Range("E14").Value = "'  -2.56%  "
$ws.Range("D15").Value = "'3.600.75"
$ws.Range("E15").Value = "'  -2.11%  "
$ws.Range("E16").Value = "'  -1.89%  "
$ws.Range("D17").Value = "'7.18"
$ws.Range("E17").Value = "'  -1.87%  "
$ws.Range("D18").Value = "'63.595.62"
$ws.Range("E18").Value = "'  -1.01%  "
$ws.Range("D19").Value = "'3.090.42"
$ws.Range("E19").Value = "'  -2.08%  "
$ws.Range("D20").Value = "'476.19"
$ws.Range("E20").Value = "'  +1.22%  "
$ws.Range("D21").Value = "'14.73"
$ws.Range("E21").Value = "'  +1.49%  "
$ws.Range("D22").Value = "'0.721"
$ws.Range("E22").Value = "'  -2.36%  "
$ws.Range("E23").Value = "'  -0.22%  "
$ws.Range("D24").Value = "'2.38"
$ws.Range("E24").Value = "'  +2.09%  "
$ws.Range("E25").Value = "'  -0.33%  "
$ws.Range("D26").Value = "'81.92"
$ws.Range("E26").Value = "'  +0.23%  "
$ws.Range("E27").Value = "'  +0.00%  "
$ws.Range("D28").Value = "'9.82"
$ws.Range("E28").Value = "'  +2.23%  "
$ws.Range("E29").Value = "'  -1.73%  "
$ws.Range("B30").Value = "'FirstDigitalUSD"
$ws.Range("C30").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "'  -0.20%  "
$ws.Range("B31").Value = "'NEARProtocol"
$ws.Range("C31").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'7.27"
$ws.Range("E31").Value = "'  -2.84%  "
$ws.Range("E32").Value = "'  -2.42%  "
$ws.Range("D33").Value = "'0.114"
$ws.Range("E33").Value = "'  +3.33%  "
$ws.Range("D34").Value = "'27.47"
$ws.Range("E34").Value = "'  -0.53%  "
$ws.Range("E35").Value = "'  +0.43%  "
$ws.Range("D36").Value = "'1.05"
$ws.Range("E36").Value = "'  -1.37%  "
$ws.Range("D37").Value = "'3.43"
$ws.Range("E37").Value = "'  +5.14%  "
$ws.Range("E38").Value = "'  -2.08%  "
$ws.Range("D39").Value = "'2.24"
$ws.Range("E39").Value = "'  -4.42%  "
$ws.Range("D40").Value = "'9.35"
$ws.Range("E40").Value = "'  +0.59%  "
$ws.Range("D41").Value = "'50.64"
$ws.Range("E41").Value = "'  -2.50%  "
$ws.Range("D42").Value = "'446.75"
$ws.Range("E42").Value = "'  -2.62%  "
$ws.Range("D43").Value = "'0.285"
$ws.Range("E43").Value = "'  -3.79%  "
$ws.Range("E44").Value = "'  -2.29%  "
$ws.Range("D45").Value = "'2.834.44"
$ws.Range("E45").Value = "'  -3.68%  "
$ws.Range("D46").Value = "'39.69"
$ws.Range("E46").Value = "'  -2.18%  "
$ws.Range("E47").Value = "'  -0.15%  "
$ws.Range("D48").Value = "'130.55"
$ws.Range("E48").Value = "'  +0.92%  "
$ws.Range("D49").Value = "'25.50"
$ws.Range("E49").Value = "'  +3.85%  "
$ws.Range("D50").Value = "'0.999"
$ws.Range("E50").Value = "'  +0.03%  "
$ws.Range("D51").Value = "'2.27"
$ws.Range("E51").Value = "'  +0.76%  "
